$d = $word.ActiveDocument

# The review comments left by "Boon Xun Sim" on the trimester/date line
# have been addressed, so remove them (and their range anchors) from the
# document entirely. Walk backwards so indices stay valid as each one is
# removed.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
